# fix: in-cell composite types
#
# Item sheet: the "attribute-1 enhancement" sub-columns (E:I) for the
# third row (item 1003) are cleared out entirely (not just emptied).
#
# Env sheet: "基本信息" becomes "基本信息1", and a new row "基本信息2" /
# "220,标题2" is inserted right after it (pushing 比例/系数 down).
#
# Selections / active tab are moved around as part of the same edit.

$wb = $excel.ActiveWorkbook

$wsItem = $wb.Worksheets.Item("Item")
$wsActivity = $wb.Worksheets.Item("Activity")
$wsEnv = $wb.Worksheets.Item("Env")

# --- Item sheet: drop the now-unused attribute-1 enhancement cells on row 5 ---
$wsItem.Range("E5:I5").Clear()

# --- Env sheet: split "基本信息" into two rows ---
$wsEnv.Range("A6").Value = "基本信息1"
$wsEnv.Rows("7").Insert()
$wsEnv.Range("A7").Value = "基本信息2"
$wsEnv.Range("B7").Value = "220,标题2"

# --- selections / active sheet ---
$wsItem.Select()
$wsItem.Range("M11").Select()

$wsActivity.Select()
$wsActivity.Range("F19").Select()

$wsEnv.Activate()
$wsEnv.Range("M3").Select()
